$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new column F "From RCSB", mirroring column E values for each data row
$ws.Range("F1").Value = "From RCSB"
$ws.Range("F2").Value = "y"
$ws.Range("F3").Value = "y"
$ws.Range("F4").Value = "n"
$ws.Range("F5").Value = "n"

$ws.Range("H7:H8").Select() | Out-Null
